$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the steel description in B2 by removing the "RME/" fragment
# from the two relevant lines (19% S/LFM+CDH/RME/H:1 -> 19% S/LFM+CDH/H:1
# and 14% S/LFM+CDM/RME/H:1 -> 14% S/LFM+CDM/H:1)
$cell = $ws.Range("B2")
$oldVal = $cell.Value()
$newVal = $oldVal.Replace("RME/", "")
$cell.Value = $newVal

# Wrap the (now multi-line) text and size the row to fit it
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.6

# Leave the active selection where the author left it after editing
$ws.Range("A10").Select() | Out-Null
